# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Ciruela - Black Amber) above the
# existing row 76, pushing the former rows 76-78 down to 78-80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 76; this shifts old rows 76,77,78
# down to 78,79,80 and keeps their formatting/styles intact.
$ws.Rows("76:77").Insert()

# --- New row 76 ---
$ws.Range("A76").Value = 9
$ws.Range("B76").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C76").Value = "Metropolitana"
$ws.Range("D76").Value = 44568
$ws.Range("E76").Value = 13
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100103
$ws.Range("H76").Value = "Frutos de hueso (carozo)"
$ws.Range("I76").Value = 100103002
$ws.Range("J76").Value = "Ciruela"
$ws.Range("K76").Value = "Black Amber"
$ws.Range("L76").Value = "Especial"
$ws.Range("M76").Value = 280
$ws.Range("N76").Value = 15000
$ws.Range("O76").Value = 15000
$ws.Range("P76").Value = 15000
$ws.Range("Q76").Value = "`$/caja 15 kilos granel"
$ws.Range("R76").Value = "Región de O'Higgins"
$ws.Range("S76").Value = 1000
$ws.Range("T76").Value = 15

# --- New row 77 ---
$ws.Range("A77").Value = 9
$ws.Range("B77").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C77").Value = "Metropolitana"
$ws.Range("D77").Value = 44568
$ws.Range("E77").Value = 13
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100103
$ws.Range("H77").Value = "Frutos de hueso (carozo)"
$ws.Range("I77").Value = 100103002
$ws.Range("J77").Value = "Ciruela"
$ws.Range("K77").Value = "Black Amber"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 350
$ws.Range("N77").Value = 12000
$ws.Range("O77").Value = 12000
$ws.Range("P77").Value = 12000
$ws.Range("Q77").Value = "`$/caja 15 kilos granel"
$ws.Range("R77").Value = "Región de O'Higgins"
$ws.Range("S77").Value = 800
$ws.Range("T77").Value = 15

"Rows inserted and populated"
